$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.5321473986446712
$ws.Range("C2").Value = 0.9894038152134554
$ws.Range("D2").Value = 0.6104588214332914
$ws.Range("G2").Value = 0.4821145882335259
$ws.Range("H2").Value = 0.992

# Row 3
$ws.Range("B3").Value = 0.2349473443729408
$ws.Range("C3").Value = 0.9954079482457484
$ws.Range("D3").Value = 0.3884558644305861
$ws.Range("G3").Value = 0.4821145882335259
$ws.Range("H3").Value = 0.992

# Row 4
$ws.Range("B4").Value = 0.2787905533223916
$ws.Range("C4").Value = 0.9946366896843218
$ws.Range("D4").Value = 0.4309113523329547
$ws.Range("G4").Value = 0.4821145882335259
$ws.Range("H4").Value = 0.992

# Row 5
$ws.Range("B5").Value = 0.4158868126993949
$ws.Range("C5").Value = 0.9917992814533819
$ws.Range("D5").Value = 0.4965755407813393
$ws.Range("G5").Value = 0.4821145882335259
$ws.Range("H5").Value = 0.992
